$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H64").Value = 6519.8
$ws.Range("I64").Value = 7000
$ws.Range("J64").Value = 4599
$ws.Range("K64").Value = 7000
$ws.Range("L64").Value = 4599
$ws.Range("M64").Value = -6752
$ws.Range("N64").Value = -5095

$ws.Range("H67").Value = 6519.8
$ws.Range("I67").Value = 7000
$ws.Range("J67").Value = 4599
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 4599
$ws.Range("M67").Value = -6142
$ws.Range("N67").Value = -6315

$ws.Range("H70").Value = 1990
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 1877.7778
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 5633.3334
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -6173.3334

$ws.Range("H73").Value = 1990
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 1877.7778
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 5633.3334
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -7505.3334

$ws.Range("H80").Value = 4346.3335
$ws.Range("I80").Value = 3436.5
$ws.Range("J80").Value = 6166
$ws.Range("K80").Value = 10309.5
$ws.Range("L80").Value = 18498
$ws.Range("M80").Value = -9311.5
$ws.Range("N80").Value = -20494

$ws.Range("H83").Value = 4346.3335
$ws.Range("I83").Value = 3436.5
$ws.Range("J83").Value = 6166
$ws.Range("K83").Value = 30928.5
$ws.Range("L83").Value = 55494
$ws.Range("M83").Value = -25936.5
$ws.Range("N83").Value = -65478

$ws.Range("H87").Value = 79998.25
$ws.Range("J87").Value = 79998.25
$ws.Range("L87").Value = 79998.25
$ws.Range("N87").Value = -82494.25

$ws.Range("H90").Value = 79998.25
$ws.Range("J90").Value = 79998.25
$ws.Range("L90").Value = 239994.75
$ws.Range("N90").Value = -252474.75

$ws.Range("H98").Value = 3687.6667
$ws.Range("I98").Value = 1741.2858
$ws.Range("J98").Value = 10500
$ws.Range("K98").Value = 1741.2858
$ws.Range("L98").Value = 10500
$ws.Range("M98").Value = -243.2858000000001
$ws.Range("N98").Value = -13496

$ws.Range("H116").Value = 7199.5
$ws.Range("I116").Value = 9900
$ws.Range("J116").Value = 4499
$ws.Range("K116").Value = 9900
$ws.Range("L116").Value = 4499
$ws.Range("M116").Value = -6458
$ws.Range("N116").Value = -11383

$ws.Range("H122").Value = 3687.6667
$ws.Range("I122").Value = 1741.2858
$ws.Range("J122").Value = 10500
$ws.Range("K122").Value = 5223.857400000001
$ws.Range("L122").Value = 31500
$ws.Range("M122").Value = -2773.857400000001
$ws.Range("N122").Value = -36400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2914.6667
$ws.Range("I45").Value = 3187.7
$ws.Range("J45").Value = 1549.5
$ws.Range("K45").Value = 3187.7
$ws.Range("L45").Value = 1549.5
$ws.Range("M45").Value = -2810.7
$ws.Range("N45").Value = -2303.5

$ws.Range("H63").Value = 2440.818
$ws.Range("I63").Value = 684.3
$ws.Range("K63").Value = 684.3
$ws.Range("M63").Value = 1.700000000000045

$ws.Range("H66").Value = 2440.818
$ws.Range("I66").Value = 684.3
$ws.Range("K66").Value = 3421.5
$ws.Range("M66").Value = 10.5

$ws.Range("H97").Value = 299.33334
$ws.Range("I97").Value = 299.5
$ws.Range("J97").Value = 299
$ws.Range("K97").Value = 299.5
$ws.Range("L97").Value = 299
$ws.Range("M97").Value = 196.5
$ws.Range("N97").Value = -1291

$ws.Range("H102").Value = 2584.1428
$ws.Range("I102").Value = 2584.1428
$ws.Range("K102").Value = 2584.1428
$ws.Range("M102").Value = -962.1428000000001

$ws.Range("H110").Value = 4166.5835
$ws.Range("J110").Value = 9166.666999999999
$ws.Range("L110").Value = 9166.666999999999
$ws.Range("N110").Value = -13256.667

$ws.Range("H122").Value = 7002.4
$ws.Range("I122").Value = 6502.75
$ws.Range("K122").Value = 19508.25
$ws.Range("M122").Value = -17058.25

$ws.Range("H139").Value = 84999.5
$ws.Range("J139").Value = 84999.5
$ws.Range("L139").Value = 84999.5
$ws.Range("N139").Value = -95279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8894.076999999999
$ws.Range("I94").Value = 7814.375
$ws.Range("J94").Value = 10621.6
$ws.Range("K94").Value = 7814.375
$ws.Range("L94").Value = 10621.6
$ws.Range("M94").Value = -7363.375
$ws.Range("N94").Value = -11523.6

$ws.Range("H107").Value = 1169.1765
$ws.Range("I107").Value = 1169.1765
$ws.Range("K107").Value = 1169.1765
$ws.Range("M107").Value = 750.8235

$ws.Range("H134").Value = 3066.6667
$ws.Range("I134").Value = 3066.6667
$ws.Range("K134").Value = 9200.000100000001
$ws.Range("M134").Value = -6665.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1556.4286
$ws.Range("J16").Value = 1474.5
$ws.Range("L16").Value = 1474.5
$ws.Range("N16").Value = -2048.5

$ws.Range("H59").Value = 7524.75
$ws.Range("I59").Value = 104
$ws.Range("J59").Value = 9998.333000000001
$ws.Range("K59").Value = 104
$ws.Range("L59").Value = 9998.333000000001
$ws.Range("M59").Value = 1041
$ws.Range("N59").Value = -12288.333

$ws.Range("H105").Value = 4535.125
$ws.Range("J105").Value = 4329.3335
$ws.Range("L105").Value = 4329.3335
$ws.Range("N105").Value = -7823.3335

$ws.Range("H113").Value = 1556.4286
$ws.Range("J113").Value = 1474.5
$ws.Range("L113").Value = 1474.5
$ws.Range("N113").Value = -5814.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17405.25
$ws.Range("I4").Value = 17405.25
$ws.Range("K4").Value = 52215.75
$ws.Range("M4").Value = -52103.75

$ws.Range("H131").Value = 2071.3
$ws.Range("J131").Value = 2676.6
$ws.Range("L131").Value = 8029.799999999999
$ws.Range("N131").Value = -18109.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11966.667
$ws.Range("I80").Value = 2950
$ws.Range("J80").Value = 30000
$ws.Range("K80").Value = 2950
$ws.Range("L80").Value = 30000
$ws.Range("M80").Value = -1952
$ws.Range("N80").Value = -31996

$ws.Range("H83").Value = 11966.667
$ws.Range("I83").Value = 2950
$ws.Range("J83").Value = 30000
$ws.Range("K83").Value = 14750
$ws.Range("L83").Value = 150000
$ws.Range("M83").Value = -9758
$ws.Range("N83").Value = -159984

$ws.Range("H113").Value = 4031.625
$ws.Range("I113").Value = 4031.625
$ws.Range("K113").Value = 4031.625
$ws.Range("M113").Value = -1861.625

$ws.Range("H140").Value = 99773.8
$ws.Range("J140").Value = 99773.8
$ws.Range("L140").Value = 99773.8
$ws.Range("N140").Value = -110133.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5730.1665
$ws.Range("I22").Value = 3844.818
$ws.Range("J22").Value = 8692.857
$ws.Range("K22").Value = 3844.818
$ws.Range("L22").Value = 8692.857
$ws.Range("M22").Value = -3549.818
$ws.Range("N22").Value = -9282.857

$ws.Range("H27").Value = 5730.1665
$ws.Range("I27").Value = 3844.818
$ws.Range("J27").Value = 8692.857
$ws.Range("K27").Value = 3844.818
$ws.Range("L27").Value = 8692.857
$ws.Range("M27").Value = -3737.818
$ws.Range("N27").Value = -8906.857

$ws.Range("H82").Value = 2800
$ws.Range("I82").Value = 2900
$ws.Range("K82").Value = 2900
$ws.Range("M82").Value = -2539

$ws.Range("H85").Value = 2800
$ws.Range("I85").Value = 2900
$ws.Range("K85").Value = 2900
$ws.Range("M85").Value = -1652

$ws.Range("H93").Value = 626.75
$ws.Range("I93").Value = 667.6667
$ws.Range("J93").Value = 504
$ws.Range("K93").Value = 667.6667
$ws.Range("L93").Value = 504
$ws.Range("M93").Value = 580.3333
$ws.Range("N93").Value = -3000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12550
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2303.75
$ws.Range("I132").Value = 1331.125
$ws.Range("K132").Value = 3993.375
$ws.Range("M132").Value = -1463.375
